# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures scraped by the scheduled runner
# to the Leviathan_Profits workbook (one hunk per Leve row, across all 8 sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 23149.5
$ws.Range("I18").Value = 6966
$ws.Range("J18").Value = 39333
$ws.Range("K18").Value = 6966
$ws.Range("L18").Value = 39333
$ws.Range("M18").Value = -6682
$ws.Range("N18").Value = -39901
# Row 40
$ws.Range("H40").Value = 2289.111
$ws.Range("I40").Value = 900.75
$ws.Range("J40").Value = 3399.8
$ws.Range("K40").Value = 900.75
$ws.Range("L40").Value = 3399.8
$ws.Range("M40").Value = -725.75
$ws.Range("N40").Value = -3749.8
# Row 82
$ws.Range("H82").Value = 484.66666
$ws.Range("I82").Value = 484.66666
$ws.Range("K82").Value = 1453.99998
$ws.Range("M82").Value = -1047.99998
# Row 85
$ws.Range("H85").Value = 484.66666
$ws.Range("I85").Value = 484.66666
$ws.Range("K85").Value = 1453.99998
$ws.Range("M85").Value = -49.99998000000005
# Row 96
$ws.Range("H96").Value = 143715.92
$ws.Range("I96").Value = 167266.17
$ws.Range("J96").Value = 2414.5
$ws.Range("K96").Value = 501798.51
$ws.Range("L96").Value = 7243.5
$ws.Range("M96").Value = -500425.51
$ws.Range("N96").Value = -9989.5
# Row 121
$ws.Range("H121").Value = 5659.6
$ws.Range("J121").Value = 5659.6
$ws.Range("L121").Value = 16978.8
$ws.Range("N121").Value = -20472.8
# Row 129
$ws.Range("H129").Value = 963.46155
$ws.Range("I129").Value = 792.5
$ws.Range("J129").Value = 1533.3334
$ws.Range("K129").Value = 2377.5
$ws.Range("L129").Value = 4600.0002
$ws.Range("M129").Value = 2622.5
$ws.Range("N129").Value = -14600.0002
# Row 131
$ws.Range("H131").Value = 297.16666
$ws.Range("I131").Value = 297.16666
$ws.Range("K131").Value = 891.4999799999999
$ws.Range("M131").Value = 4148.50002
# Row 138
$ws.Range("H138").Value = 2192.8
$ws.Range("I138").Value = 1314.6666
$ws.Range("J138").Value = 4827.2
$ws.Range("K138").Value = 3943.9998
$ws.Range("L138").Value = 14481.6
$ws.Range("M138").Value = 1196.0002
$ws.Range("N138").Value = -24761.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 128734.63
$ws.Range("I32").Value = 134733.64
$ws.Range("K32").Value = 134733.64
$ws.Range("M32").Value = -134446.64
# Row 74
$ws.Range("H74").Value = 1105.3549
$ws.Range("I74").Value = 911.375
$ws.Range("J74").Value = 1770.4286
$ws.Range("K74").Value = 911.375
$ws.Range("L74").Value = 1770.4286
$ws.Range("M74").Value = -37.375
$ws.Range("N74").Value = -3518.4286
# Row 77
$ws.Range("H77").Value = 1105.3549
$ws.Range("I77").Value = 911.375
$ws.Range("J77").Value = 1770.4286
$ws.Range("K77").Value = 4556.875
$ws.Range("L77").Value = 8852.143
$ws.Range("M77").Value = -188.875
$ws.Range("N77").Value = -17588.143
# Row 111
$ws.Range("H111").Value = 64762
$ws.Range("J111").Value = 64762
$ws.Range("L111").Value = 64762
$ws.Range("N111").Value = -72942
# Row 119
$ws.Range("H119").Value = 35931.668
$ws.Range("J119").Value = 35931.668
$ws.Range("L119").Value = 35931.668
$ws.Range("N119").Value = -45607.668

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 3608.5454
$ws.Range("I22").Value = 3608.5454
$ws.Range("K22").Value = 3608.5454
$ws.Range("M22").Value = -3435.5454
# Row 94
$ws.Range("H94").Value = 1090.35
$ws.Range("I94").Value = 861.25
$ws.Range("K94").Value = 861.25
$ws.Range("M94").Value = -410.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2399.6
$ws.Range("I31").Value = 1926.5333
$ws.Range("K31").Value = 1926.5333
$ws.Range("M31").Value = -1631.5333
# Row 34
$ws.Range("H34").Value = 2399.6
$ws.Range("I34").Value = 1926.5333
$ws.Range("K34").Value = 1926.5333
$ws.Range("M34").Value = -1724.5333
# Row 38
$ws.Range("H38").Value = 8412.546
$ws.Range("I38").Value = 5519
$ws.Range("K38").Value = 5519
$ws.Range("M38").Value = -5142
# Row 46
$ws.Range("H46").Value = 8412.546
$ws.Range("I46").Value = 5519
$ws.Range("K46").Value = 5519
$ws.Range("M46").Value = -5308

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 384.66666
$ws.Range("I7").Value = 337.5
$ws.Range("K7").Value = 1012.5
$ws.Range("M7").Value = -900.5
# Row 23
$ws.Range("H23").Value = 716.6923
$ws.Range("I23").Value = 548
$ws.Range("J23").Value = 767.3
$ws.Range("K23").Value = 1644
$ws.Range("L23").Value = 2301.9
$ws.Range("M23").Value = -1409
$ws.Range("N23").Value = -2771.9
# Row 107
$ws.Range("H107").Value = 884.5
$ws.Range("I107").Value = 1116.6666
$ws.Range("J107").Value = 821.1818
$ws.Range("K107").Value = 3349.9998
$ws.Range("L107").Value = 2463.5454
$ws.Range("M107").Value = -1429.9998
$ws.Range("N107").Value = -6303.5454
# Row 131
$ws.Range("H131").Value = 6007.7827
$ws.Range("I131").Value = 10268.272
$ws.Range("J131").Value = 2102.3333
$ws.Range("K131").Value = 30804.816
$ws.Range("L131").Value = 6306.999899999999
$ws.Range("M131").Value = -25764.816
$ws.Range("N131").Value = -16386.9999
# Row 132
$ws.Range("H132").Value = 2625
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2625
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 23625
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -28685

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 27077
$ws.Range("I97").Value = 35226.473
$ws.Range("J97").Value = 1270.3334
$ws.Range("K97").Value = 35226.473
$ws.Range("L97").Value = 1270.3334
$ws.Range("M97").Value = -34730.473
$ws.Range("N97").Value = -2262.3334

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 55540
$ws.Range("I46").Value = 72553.336
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 72553.336
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -72365.336
$ws.Range("N46").Value = -4876
# Row 55
$ws.Range("H55").Value = 322.36667
$ws.Range("I55").Value = 241.1875
$ws.Range("J55").Value = 415.14285
$ws.Range("K55").Value = 241.1875
$ws.Range("L55").Value = 415.14285
$ws.Range("M55").Value = -68.1875
$ws.Range("N55").Value = -761.14285
# Row 93
$ws.Range("H93").Value = 69026.2
$ws.Range("I93").Value = 2949.5
$ws.Range("K93").Value = 2949.5
$ws.Range("M93").Value = -1701.5

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 40891.75
$ws.Range("I45").Value = 37784.5
$ws.Range("J45").Value = 43999
$ws.Range("K45").Value = 37784.5
$ws.Range("L45").Value = 43999
$ws.Range("M45").Value = -37293.5
$ws.Range("N45").Value = -44981
# Row 115
$ws.Range("H115").Value = 80000
$ws.Range("J115").Value = 80000
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -83134
# Row 132
$ws.Range("H132").Value = 5893.5674
$ws.Range("I132").Value = 7948.5
$ws.Range("K132").Value = 23845.5
$ws.Range("M132").Value = -21315.5
